# Contest 8 PBKS vs CSK — enter each player's raw points for match row 20.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E20").Value = 40    # Jaya
$ws.Range("H20").Value = 80    # Justin
$ws.Range("K20").Value = 0     # Ram
$ws.Range("N20").Value = 20    # Sibi
$ws.Range("Q20").Value = 60    # Sundar
$ws.Range("T20").Value = 30    # Balaji
$ws.Range("W20").Value = 70    # Upili
$ws.Range("Z20").Value = 100   # Vicky
$ws.Range("AC20").Value = 50   # Raghu
